$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 900
$ws.Cells.Item(9, 9).Value = 900
$ws.Cells.Item(9, 10).Value = 900
$ws.Cells.Item(9, 11).Value = 900
$ws.Cells.Item(9, 12).Value = 900
$ws.Cells.Item(9, 13).Value = -731
$ws.Cells.Item(9, 14).Value = -1238

$ws.Cells.Item(28, 8).Value = 1304.25
$ws.Cells.Item(28, 9).Value = 449.16666
$ws.Cells.Item(28, 11).Value = 449.16666
$ws.Cells.Item(28, 13).Value = 35.83334000000002

$ws.Cells.Item(33, 8).Value = 172.6
$ws.Cells.Item(33, 9).Value = 173.94737
$ws.Cells.Item(33, 11).Value = 173.94737
$ws.Cells.Item(33, 13).Value = 55.05262999999999

$ws.Cells.Item(76, 8).Value = 7999
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 7999
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 7999
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(76, 14).Value = -8629

$ws.Cells.Item(79, 8).Value = 7999
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 7999
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 7999
$ws.Cells.Item(79, 13).ClearContents()
$ws.Cells.Item(79, 14).Value = -10183

$ws.Cells.Item(113, 8).Value = 2704.318
$ws.Cells.Item(113, 9).Value = 3266.3333
$ws.Cells.Item(113, 10).Value = 2615.5789
$ws.Cells.Item(113, 11).Value = 3266.3333
$ws.Cells.Item(113, 12).Value = 2615.5789
$ws.Cells.Item(113, 13).Value = -12.33329999999978
$ws.Cells.Item(113, 14).Value = -9123.5789

$ws.Cells.Item(137, 8).Value = 2389340.8
$ws.Cells.Item(137, 9).Value = 8334383.5
$ws.Cells.Item(137, 11).Value = 25003150.5
$ws.Cells.Item(137, 13).Value = -25000600.5

$ws.Cells.Item(138, 8).Value = 3119.1775
$ws.Cells.Item(138, 9).Value = 4368
$ws.Cells.Item(138, 10).Value = 2684.8044
$ws.Cells.Item(138, 11).Value = 13104
$ws.Cells.Item(138, 12).Value = 8054.4132
$ws.Cells.Item(138, 13).Value = -7964
$ws.Cells.Item(138, 14).Value = -18334.4132

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5453.515
$ws.Cells.Item(32, 9).Value = 5010.8696
$ws.Cells.Item(32, 10).Value = 6471.6
$ws.Cells.Item(32, 11).Value = 5010.8696
$ws.Cells.Item(32, 12).Value = 6471.6
$ws.Cells.Item(32, 13).Value = -4723.8696
$ws.Cells.Item(32, 14).Value = -7045.6

$ws.Cells.Item(45, 8).Value = 36957.5
$ws.Cells.Item(45, 9).Value = 39770.547
$ws.Cells.Item(45, 11).Value = 39770.547
$ws.Cells.Item(45, 13).Value = -39393.547

$ws.Cells.Item(61, 8).Value = 2793.6553
$ws.Cells.Item(61, 9).Value = 1850.1
$ws.Cells.Item(61, 11).Value = 1850.1
$ws.Cells.Item(61, 13).Value = -1638.1

$ws.Cells.Item(74, 8).Value = 113286.66
$ws.Cells.Item(74, 9).Value = 186262.77
$ws.Cells.Item(74, 11).Value = 186262.77
$ws.Cells.Item(74, 13).Value = -185388.77

$ws.Cells.Item(77, 8).Value = 113286.66
$ws.Cells.Item(77, 9).Value = 186262.77
$ws.Cells.Item(77, 11).Value = 931313.85
$ws.Cells.Item(77, 13).Value = -926945.85

$ws.Cells.Item(96, 8).Value = 44999.5
$ws.Cells.Item(96, 10).Value = 44999.5
$ws.Cells.Item(96, 12).Value = 44999.5
$ws.Cells.Item(96, 14).Value = -50491.5

$ws.Cells.Item(132, 8).Value = 2890.037
$ws.Cells.Item(132, 9).Value = 1976.55
$ws.Cells.Item(132, 10).Value = 5500
$ws.Cells.Item(132, 11).Value = 5929.65
$ws.Cells.Item(132, 12).Value = 16500
$ws.Cells.Item(132, 13).Value = -3399.65
$ws.Cells.Item(132, 14).Value = -21560

$ws.Cells.Item(136, 8).Value = 2793.6553
$ws.Cells.Item(136, 9).Value = 1850.1
$ws.Cells.Item(136, 11).Value = 5550.299999999999
$ws.Cells.Item(136, 13).Value = -3000.299999999999

$ws.Cells.Item(139, 8).Value = 84529.82000000001
$ws.Cells.Item(139, 10).Value = 84529.82000000001
$ws.Cells.Item(139, 12).Value = 84529.82000000001
$ws.Cells.Item(139, 14).Value = -94809.82000000001

$ws.Cells.Item(140, 8).Value = 77854.86
$ws.Cells.Item(140, 10).Value = 85830.664
$ws.Cells.Item(140, 12).Value = 85830.664
$ws.Cells.Item(140, 14).Value = -96190.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1520.9642
$ws.Cells.Item(107, 9).Value = 1306.5238
$ws.Cells.Item(107, 10).Value = 2164.2856
$ws.Cells.Item(107, 11).Value = 1306.5238
$ws.Cells.Item(107, 12).Value = 2164.2856
$ws.Cells.Item(107, 13).Value = 613.4762000000001
$ws.Cells.Item(107, 14).Value = -6004.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4437.6587
$ws.Cells.Item(31, 9).Value = 2878.1667
$ws.Cells.Item(31, 10).Value = 6639.294
$ws.Cells.Item(31, 11).Value = 2878.1667
$ws.Cells.Item(31, 12).Value = 6639.294
$ws.Cells.Item(31, 13).Value = -2583.1667
$ws.Cells.Item(31, 14).Value = -7229.294

$ws.Cells.Item(34, 8).Value = 4437.6587
$ws.Cells.Item(34, 9).Value = 2878.1667
$ws.Cells.Item(34, 10).Value = 6639.294
$ws.Cells.Item(34, 11).Value = 2878.1667
$ws.Cells.Item(34, 12).Value = 6639.294
$ws.Cells.Item(34, 13).Value = -2676.1667
$ws.Cells.Item(34, 14).Value = -7043.294

$ws.Cells.Item(94, 8).Value = 2169.2
$ws.Cells.Item(94, 9).Value = 1909.2
$ws.Cells.Item(94, 11).Value = 1909.2
$ws.Cells.Item(94, 13).Value = -1458.2

$ws.Cells.Item(99, 8).Value = 5083.8887
$ws.Cells.Item(99, 9).Value = 3826.1667
$ws.Cells.Item(99, 11).Value = 3826.1667
$ws.Cells.Item(99, 13).Value = -2328.1667

$ws.Cells.Item(126, 8).Value = 5083.8887
$ws.Cells.Item(126, 9).Value = 3826.1667
$ws.Cells.Item(126, 11).Value = 11478.5001
$ws.Cells.Item(126, 13).Value = -9008.500100000001

$ws.Cells.Item(132, 8).Value = 2394.1
$ws.Cells.Item(132, 9).Value = 1722.1538
$ws.Cells.Item(132, 10).Value = 3642
$ws.Cells.Item(132, 11).Value = 5166.4614
$ws.Cells.Item(132, 12).Value = 10926
$ws.Cells.Item(132, 13).Value = -2636.4614
$ws.Cells.Item(132, 14).Value = -15986

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1767553.5
$ws.Cells.Item(4, 9).Value = 1152393.6
$ws.Cells.Item(4, 10).Value = 7303992
$ws.Cells.Item(4, 11).Value = 3457180.8
$ws.Cells.Item(4, 12).Value = 21911976
$ws.Cells.Item(4, 13).Value = -3457068.8
$ws.Cells.Item(4, 14).Value = -21912200

$ws.Cells.Item(44, 8).Value = 2244.75
$ws.Cells.Item(44, 10).Value = 2244.75
$ws.Cells.Item(44, 12).Value = 6734.25
$ws.Cells.Item(44, 14).Value = -7530.25

$ws.Cells.Item(95, 8).Value = 9996
$ws.Cells.Item(95, 10).Value = 9996
$ws.Cells.Item(95, 12).Value = 29988
$ws.Cells.Item(95, 14).Value = -34106

$ws.Cells.Item(100, 8).Value = 3008.3333
$ws.Cells.Item(100, 10).Value = 3008.3333
$ws.Cells.Item(100, 12).Value = 9024.999899999999
$ws.Cells.Item(100, 14).Value = -10646.9999

$ws.Cells.Item(129, 8).Value = 81077.57000000001
$ws.Cells.Item(129, 9).Value = 424.5
$ws.Cells.Item(129, 10).Value = 113338.8
$ws.Cells.Item(129, 11).Value = 1273.5
$ws.Cells.Item(129, 12).Value = 340016.4
$ws.Cells.Item(129, 13).Value = 3726.5
$ws.Cells.Item(129, 14).Value = -350016.4

$ws.Cells.Item(134, 8).Value = 798.3333
$ws.Cells.Item(134, 9).Value = 798.3333
$ws.Cells.Item(134, 11).Value = 2394.9999
$ws.Cells.Item(134, 13).Value = 2675.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(96, 8).Value = 56340.332
$ws.Cells.Item(96, 10).Value = 56340.332
$ws.Cells.Item(96, 12).Value = 56340.332
$ws.Cells.Item(96, 14).Value = -61832.332

$ws.Cells.Item(122, 8).Value = 5102.476
$ws.Cells.Item(122, 9).Value = 3301.1428
$ws.Cells.Item(122, 11).Value = 9903.428400000001
$ws.Cells.Item(122, 13).Value = -7453.428400000001

$ws.Cells.Item(132, 8).Value = 3601.2273
$ws.Cells.Item(132, 9).Value = 3330.5715
$ws.Cells.Item(132, 11).Value = 9991.7145
$ws.Cells.Item(132, 13).Value = -7461.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 771.2857
$ws.Cells.Item(22, 10).Value = 799.8
$ws.Cells.Item(22, 12).Value = 799.8
$ws.Cells.Item(22, 14).Value = -1389.8

$ws.Cells.Item(27, 8).Value = 771.2857
$ws.Cells.Item(27, 10).Value = 799.8
$ws.Cells.Item(27, 12).Value = 799.8
$ws.Cells.Item(27, 14).Value = -1013.8

$ws.Cells.Item(61, 8).Value = 9166.308000000001
$ws.Cells.Item(61, 9).Value = 784.7273
$ws.Cells.Item(61, 11).Value = 784.7273
$ws.Cells.Item(61, 13).Value = -582.7273

$ws.Cells.Item(113, 8).Value = 9166.308000000001
$ws.Cells.Item(113, 9).Value = 784.7273
$ws.Cells.Item(113, 11).Value = 784.7273
$ws.Cells.Item(113, 13).Value = 1385.2727

$ws.Cells.Item(132, 8).Value = 6541.8423
$ws.Cells.Item(132, 9).Value = 2732
$ws.Cells.Item(132, 11).Value = 8196
$ws.Cells.Item(132, 13).Value = -5666

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1707.4
$ws.Cells.Item(107, 9).Value = 929
$ws.Cells.Item(107, 10).Value = 2875
$ws.Cells.Item(107, 11).Value = 2787
$ws.Cells.Item(107, 12).Value = 8625
$ws.Cells.Item(107, 13).Value = -867
$ws.Cells.Item(107, 14).Value = -12465
